# Replace Product, IT, and Finance templates with correct industry-specific content
# Updates the IT Training Schedule workbook: renames AI/ML-specific module and
# audience labels to generic IT equivalents across the Overview and Detailed
# Training Schedule sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Training Schedule Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Schedule Overview")

# Touch row 6 and row 17 (blank spacer rows) so they are materialized.
$ws1.Rows.Item(6).Font.Bold = $false
$ws1.Rows.Item(17).Font.Bold = $false

$ws1.Range("A9").Value  = "IT Fundamentals (AI-101)"
$ws1.Range("A10").Value = "IT Platform Overview (AI-102)"

$ws1.Range("B11").Value = "System Administrators"
$ws1.Range("B12").Value = "IT Managers"
$ws1.Range("B13").Value = "DevOps Engineers, IT"
$ws1.Range("B14").Value = "DevOps Engineers, QA"

# ---------------------------------------------------------------------------
# Sheet 2: "Detailed Training Schedule"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Detailed Training Schedule")

# Touch row 2 (blank spacer row) so it is materialized.
$ws2.Rows.Item(2).Font.Bold = $false

$ws2.Range("B4").Value = "IT Fundamentals"
$ws2.Range("B5").Value = "IT Platform Overview"

$ws2.Range("C6").Value  = "System Administrators"
$ws2.Range("C7").Value  = "System Administrators"
$ws2.Range("C8").Value  = "System Administrators"
$ws2.Range("C9").Value  = "IT Managers"
$ws2.Range("C10").Value = "IT Managers"
$ws2.Range("C11").Value = "DevOps Engineers, IT"
$ws2.Range("C12").Value = "DevOps Engineers, IT"
$ws2.Range("C13").Value = "DevOps Engineers, QA"
$ws2.Range("C14").Value = "DevOps Engineers, QA"

# ---------------------------------------------------------------------------
# Sheet 3: "Instructor Schedule" - touch row 2 (blank spacer row)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Instructor Schedule")
$ws3.Rows.Item(2).Font.Bold = $false

# ---------------------------------------------------------------------------
# Sheet 4: "Facility Schedule" - touch row 2 (blank spacer row)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Facility Schedule")
$ws4.Rows.Item(2).Font.Bold = $false

# ---------------------------------------------------------------------------
# Sheet 5: "Participant Tracking" - touch row 2 (blank spacer row)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Participant Tracking")
$ws5.Rows.Item(2).Font.Bold = $false

Write-Host "Edit complete"
